{"js": "// Update the date line and the 20x5 arithmetic-answer table to the new\n// values described by the diff. Positional (ordinal) replacement is used\n// throughout because several old strings repeat (e.g. \"91-38=53\" occurs\n// twice) yet map to different new values depending on position, so a\n// global text find/replace would be unsafe.\n\nconst body = context.document.body;\n\n// 1) First paragraph: the date line \"2023-11-28 Tuesday\" -> \"2023-11-29 Wednesday\"\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.load(\"text\");\nawait context.sync();\n\nif (firstParagraph.text.trim() === \"2023-11-28 Tuesday\") {\n  firstParagraph.insertText(\"2023-11-29 Wednesday\", Word.InsertLocation.replace);\n}\n\n// 2) The table of arithmetic problems: replace all 100 cell values in place,\n//    preserving row/column position (this is what \"values\" addresses).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = [\n  [\"74-25=49\", \"18+73=91\", \"22-7=15\", \"90-35=55\", \"18+15=33\"],\n  [\"41-12=29\", \"15+76=91\", \"13+78=91\", \"93-76=17\", \"22+9=31\"],\n  [\"70-4=66\", \"56+17=73\", \"54-7=47\", \"17+66=83\", \"49+14=63\"],\n  [\"36+17=53\", \"5+9=14\", \"7+69=76\", \"62-19=43\", \"49+35=84\"],\n  [\"7+8=15\", \"75+18=93\", \"74-5=69\", \"92-17=75\", \"30-19=11\"],\n  [\"52-33=19\", \"85-46=39\", \"40-2=38\", \"5+69=74\", \"82-53=29\"],\n  [\"7+84=91\", \"11-9=2\", \"62-14=48\", \"62-17=45\", \"16+48=64\"],\n  [\"82-35=47\", \"91-18=73\", \"29+39=68\", \"63+19=82\", \"40-6=34\"],\n  [\"25-19=6\", \"23+58=81\", \"7+58=65\", \"45-9=36\", \"8+15=23\"],\n  [\"91-68=23\", \"13-9=4\", \"41-32=9\", \"30-21=9\", \"91-52=39\"],\n  [\"63-17=46\", \"8+46=54\", \"46+19=65\", \"90-4=86\", \"83-55=28\"],\n  [\"60-18=42\", \"8+56=64\", \"82-69=13\", \"20-13=7\", \"92-45=47\"],\n  [\"63-9=54\", \"76-57=19\", \"90-71=19\", \"85-16=69\", \"10-5=5\"],\n  [\"93-78=15\", \"88-59=29\", \"82-37=45\", \"75+16=91\", \"81-44=37\"],\n  [\"98-89=9\", \"47-19=28\", \"82-56=26\", \"81-3=78\", \"57+15=72\"],\n  [\"23+18=41\", \"81-73=8\", \"22+29=51\", \"58+35=93\", \"92-8=84\"],\n  [\"15+28=43\", \"63-6=57\", \"24-18=6\", \"18+25=43\", \"9+58=67\"],\n  [\"95-17=78\", \"61-24=37\", \"72-16=56\", \"49+34=83\", \"74-65=9\"],\n  [\"53-28=25\", \"63-7=56\", \"64-38=26\", \"62-33=29\", \"94-45=49\"],\n  [\"92-87=5\", \"25+67=92\", \"27+49=76\", \"73-46=27\", \"18+79=97\"],\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the date line and the 20x5 arithmetic-answer table to the new\n# values described by the diff. Cells are addressed by explicit (row,\n# column) position because several old strings repeat (e.g. \"91-38=53\"\n# occurs twice) yet map to different new values depending on position, so\n# a global text find/replace would be unsafe.\n\n$d = $word.ActiveDocument\n\n# 1) First paragraph: the date line \"2023-11-28 Tuesday\" -> \"2023-11-29 Wednesday\"\n$firstParagraph = $d.Paragraphs.Item(1)\nif ($firstParagraph.Range.Text.Trim() -eq \"2023-11-28 Tuesday\") {\n    $firstParagraph.Range.Text = \"2023-11-29 Wednesday\"\n}\n\n# 2) The table of arithmetic problems: replace all 100 cell values in place,\n#    preserving row/column position (this is what Cell(row, col) addresses).\n$newValues = @(\n    @(\"74-25=49\", \"18+73=91\", \"22-7=15\", \"90-35=55\", \"18+15=33\"),\n    @(\"41-12=29\", \"15+76=91\", \"13+78=91\", \"93-76=17\", \"22+9=31\"),\n    @(\"70-4=66\", \"56+17=73\", \"54-7=47\", \"17+66=83\", \"49+14=63\"),\n    @(\"36+17=53\", \"5+9=14\", \"7+69=76\", \"62-19=43\", \"49+35=84\"),\n    @(\"7+8=15\", \"75+18=93\", \"74-5=69\", \"92-17=75\", \"30-19=11\"),\n    @(\"52-33=19\", \"85-46=39\", \"40-2=38\", \"5+69=74\", \"82-53=29\"),\n    @(\"7+84=91\", \"11-9=2\", \"62-14=48\", \"62-17=45\", \"16+48=64\"),\n    @(\"82-35=47\", \"91-18=73\", \"29+39=68\", \"63+19=82\", \"40-6=34\"),\n    @(\"25-19=6\", \"23+58=81\", \"7+58=65\", \"45-9=36\", \"8+15=23\"),\n    @(\"91-68=23\", \"13-9=4\", \"41-32=9\", \"30-21=9\", \"91-52=39\"),\n    @(\"63-17=46\", \"8+46=54\", \"46+19=65\", \"90-4=86\", \"83-55=28\"),\n    @(\"60-18=42\", \"8+56=64\", \"82-69=13\", \"20-13=7\", \"92-45=47\"),\n    @(\"63-9=54\", \"76-57=19\", \"90-71=19\", \"85-16=69\", \"10-5=5\"),\n    @(\"93-78=15\", \"88-59=29\", \"82-37=45\", \"75+16=91\", \"81-44=37\"),\n    @(\"98-89=9\", \"47-19=28\", \"82-56=26\", \"81-3=78\", \"57+15=72\"),\n    @(\"23+18=41\", \"81-73=8\", \"22+29=51\", \"58+35=93\", \"92-8=84\"),\n    @(\"15+28=43\", \"63-6=57\", \"24-18=6\", \"18+25=43\", \"9+58=67\"),\n    @(\"95-17=78\", \"61-24=37\", \"72-16=56\", \"49+34=83\", \"74-65=9\"),\n    @(\"53-28=25\", \"63-7=56\", \"64-38=26\", \"62-33=29\", \"94-45=49\"),\n    @(\"92-87=5\", \"25+67=92\", \"27+49=76\", \"73-46=27\", \"18+79=97\")\n)\n\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le 20; $r++) {\n    for ($c = 1; $c -le 5; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
